{"js": "const pairs = [\n  [\"2025-04-03 Thursday\", \"2025-04-04 Friday\"],\n  [\"782\u00d78=6256\", \"586\u00d76=3516\"],\n  [\"360\u00d72=720\", \"257\u00d78=2056\"],\n  [\"475\u00d74=1900\", \"336\u00d73=1008\"],\n  [\"720\u00d77=5040\", \"808\u00d77=5656\"],\n  [\"228\u00d72=456\", \"736\u00d79=6624\"],\n  [\"403\u00d78=3224\", \"428\u00d74=1712\"],\n  [\"543\u00d75=2715\", \"954\u00d79=8586\"],\n  [\"946\u00d73=2838\", \"115\u00d76=690\"],\n  [\"991\u00d77=6937\", \"216\u00d72=432\"],\n  [\"204\u00d76=1224\", \"302\u00d77=2114\"],\n  [\"799\u00d78=6392\", \"458\u00d75=2290\"],\n  [\"370\u00d78=2960\", \"589\u00d76=3534\"],\n  [\"453\u00d73=1359\", \"289\u00d78=2312\"],\n  [\"982\u00d79=8838\", \"833\u00d76=4998\"],\n  [\"857\u00d77=5999\", \"374\u00d79=3366\"],\n  [\"438\u00d74=1752\", \"265\u00d72=530\"],\n  [\"568\u00d77=3976\", \"913\u00d73=2739\"],\n  [\"726\u00d78=5808\", \"421\u00d73=1263\"],\n  [\"911\u00d75=4555\", \"395\u00d79=3555\"],\n  [\"735\u00d76=4410\", \"771\u00d75=3855\"],\n  [\"617\u00d78=4936\", \"251\u00d79=2259\"],\n  [\"236\u00d79=2124\", \"558\u00d78=4464\"],\n  [\"169\u00d77=1183\", \"544\u00d73=1632\"],\n  [\"162\u00d72=324\", \"166\u00d77=1162\"],\n  [\"734\u00d78=5872\", \"964\u00d74=3856\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, 'Replace');\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2025-04-03 Thursday\", \"2025-04-04 Friday\")\n    ,@(\"782\u00d78=6256\", \"586\u00d76=3516\")\n    ,@(\"360\u00d72=720\", \"257\u00d78=2056\")\n    ,@(\"475\u00d74=1900\", \"336\u00d73=1008\")\n    ,@(\"720\u00d77=5040\", \"808\u00d77=5656\")\n    ,@(\"228\u00d72=456\", \"736\u00d79=6624\")\n    ,@(\"403\u00d78=3224\", \"428\u00d74=1712\")\n    ,@(\"543\u00d75=2715\", \"954\u00d79=8586\")\n    ,@(\"946\u00d73=2838\", \"115\u00d76=690\")\n    ,@(\"991\u00d77=6937\", \"216\u00d72=432\")\n    ,@(\"204\u00d76=1224\", \"302\u00d77=2114\")\n    ,@(\"799\u00d78=6392\", \"458\u00d75=2290\")\n    ,@(\"370\u00d78=2960\", \"589\u00d76=3534\")\n    ,@(\"453\u00d73=1359\", \"289\u00d78=2312\")\n    ,@(\"982\u00d79=8838\", \"833\u00d76=4998\")\n    ,@(\"857\u00d77=5999\", \"374\u00d79=3366\")\n    ,@(\"438\u00d74=1752\", \"265\u00d72=530\")\n    ,@(\"568\u00d77=3976\", \"913\u00d73=2739\")\n    ,@(\"726\u00d78=5808\", \"421\u00d73=1263\")\n    ,@(\"911\u00d75=4555\", \"395\u00d79=3555\")\n    ,@(\"735\u00d76=4410\", \"771\u00d75=3855\")\n    ,@(\"617\u00d78=4936\", \"251\u00d79=2259\")\n    ,@(\"236\u00d79=2124\", \"558\u00d78=4464\")\n    ,@(\"169\u00d77=1183\", \"544\u00d73=1632\")\n    ,@(\"162\u00d72=324\", \"166\u00d77=1162\")\n    ,@(\"734\u00d78=5872\", \"964\u00d74=3856\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Find/Replace failed for: $oldText\"\n    }\n}\n"}
